$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells (Coin, Link, Price, Volume(1h)) to match latest scrape.
# Column D (Price) values are force-prefixed with a leading apostrophe so that
# Excel stores them as text (preserving formats like "1.00" / "18.30") instead
# of auto-converting them to numbers.

$ws.Range("D2").Value = "'56.761.59"
$ws.Range("E2").Value = "  +3.34%  "
$ws.Range("D3").Value = "'2.327.04"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'521.50"
$ws.Range("E5").Value = "  +2.80%  "
$ws.Range("D6").Value = "'134.99"
$ws.Range("E6").Value = "  +4.02%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'2.352.64"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("D13").Value = "'0.342"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'23.83"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "'2.766.45"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "'56.832.84"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "'2.342.48"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'10.49"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("D20").Value = "'4.23"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "'323.77"
$ws.Range("E21").Value = "  +3.97%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "'60.65"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'0.164"
$ws.Range("E25").Value = "  +8.12%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'7.91"
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("D28").Value = "'1.29"
$ws.Range("E28").Value = "  +12.47%  "
$ws.Range("D29").Value = "'0.0₃0752"
$ws.Range("E29").Value = "  +5.91%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.73"
$ws.Range("E30").Value = "  +5.84%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'170.14"
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").Value = "'6.18"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'18.30"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'0.991"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").Value = "'0.924"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'4.04"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +7.89%  "
$ws.Range("D40").Value = "'37.97"
$ws.Range("E40").Value = "  +3.17%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'3.59"
$ws.Range("E42").Value = "  +4.59%  "
$ws.Range("D43").Value = "'138.14"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'279.19"
$ws.Range("E44").Value = "  +7.37%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'5.18"
$ws.Range("E45").Value = "  +5.36%  "
$ws.Range("D46").Value = "'0.0934"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").Value = "'0.0505"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "'0.563"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").Value = "'0.0218"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").Value = "'17.85"
$ws.Range("E50").Value = "  +7.67%  "
$ws.Range("E51").Value = "  +0.23%  "
